$wb = $excel.ActiveWorkbook

# The active sheet is "Sheet3" (3rd tab), which matches the diff context.
$ws = $wb.Worksheets.Item("Sheet3")

# Fix the encoding/decoding of the html code in C1: was shared string "arabic",
# now becomes the literal text "&#123;&#55;"
$ws.Range("C1").Value = "&#123;&#55;"

# Update the active selection on that sheet to G5
$ws.Activate()
$ws.Range("G5").Select()
